# This script applies a swap of species-record data between rows 14/16
# and between rows 17/18 on the "Artfynd" sheet, matching the target
# OOXML diff. Only the cells that actually differ between the two rows
# in each pair are touched; shared fields (location, dates, etc.) are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 14 <-> Row 16 swap ----
$ws.Range("A14").Value = 111798755
$ws.Range("B14").Value = 90709
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 5448
$ws.Range("F14").Value = "Svartvit taggsvamp"
$ws.Range("G14").Value = "Phellodon connatus"
$ws.Range("H14").Value = "(Schultz) nom.prov"
$ws.Range("Q14").Value = 753030.7189070459
$ws.Range("R14").Value = 7090920.781295684
$ws.Range("S14").Value = 25
$ws.Range("AF14").Value = ""
$ws.Range("AI14").Value = ""

$ws.Range("A16").Value = 111798757
$ws.Range("B16").Value = 81076
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 5046
$ws.Range("F16").Value = "Grön jordtunga"
$ws.Range("G16").Value = "Microglossum viride"
$ws.Range("H16").Value = "(Pers.:Fr.) Gillet"
$ws.Range("Q16").Value = 753108.8301749222
$ws.Range("R16").Value = 7091007.708399305
$ws.Range("S16").Value = 100
$ws.Range("AF16").Value = "mikroskoperad"
$ws.Range("AI16").Value = "Granskog"

# ---- Row 17 <-> Row 18 swap ----
$ws.Range("A17").Value = 111961716
$ws.Range("B17").Value = 81076
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 5046
$ws.Range("F17").Value = "Grön jordtunga"
$ws.Range("G17").Value = "Microglossum viride"
$ws.Range("H17").Value = "(Pers.:Fr.) Gillet"
# "Antal" (count) is stored as text throughout this sheet, not a number -
# force text formatting so "2" round-trips as a string, not 2 (number).
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "2"
$ws.Range("J17").Value = "mycel"
$ws.Range("AF17").Value = "mikroskoperad"
$ws.Range("AX17").Value = "Stefan Phalagorn Bergström, Andreas Estensen, Annika  Carlberg , Ola Elleström, Thomas Strid, Anne Järvinen, Emma Sewell"

$ws.Range("A18").Value = 111961472
$ws.Range("B18").Value = 90709
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 5448
$ws.Range("F18").Value = "Svartvit taggsvamp"
$ws.Range("G18").Value = "Phellodon connatus"
$ws.Range("H18").Value = "(Schultz) nom.prov"
$ws.Range("I18").Value = ""
$ws.Range("J18").Value = ""
$ws.Range("AF18").Value = ""
$ws.Range("AX18").Value = "Stefan Phalagorn Bergström, Annika  Carlberg , Andreas Estensen, Ola Elleström, Anne Järvinen, Emma Sewell, Thomas Strid"
